# Weekly update: insert a new week's worth of Acelga price rows (Extra /
# Primera / Segunda) at the top of the existing data block, pushing the
# rest of the historical rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right before row 360 (old rows 360:383 move to 363:386)
$ws.Rows.Item(360).Resize(3).EntireRow.Insert()

# Use the (now shifted) row 363 - which holds the same static template values
# (Mercado, Region, Codreg, Categoria, Variedad, Unidad, Origen, Kg o Unidades,
# Clasificacion) shared by every row in this block - as a formatting/content
# template for the 3 new rows.
$ws.Range("A363:R363").Copy()
$ws.Range("A360:R360").PasteSpecial(-4104)
$ws.Range("A363:R363").Copy()
$ws.Range("A361:R361").PasteSpecial(-4104)
$ws.Range("A363:R363").Copy()
$ws.Range("A362:R362").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Row 360: Acelga, Extra
$ws.Range("D360").Value = 44516
$ws.Range("I360").Value = "Extra"
$ws.Range("J360").Value = 25
$ws.Range("K360").Value = 12000
$ws.Range("L360").Value = 12000
$ws.Range("M360").Value = 12000
$ws.Range("P360").Value = 4000

# Row 361: Acelga, Primera
$ws.Range("D361").Value = 44516
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 79
$ws.Range("K361").Value = 10000
$ws.Range("L361").Value = 11000
$ws.Range("M361").Value = 10494
$ws.Range("P361").Value = 3498

# Row 362: Acelga, Segunda
$ws.Range("D362").Value = 44516
$ws.Range("I362").Value = "Segunda"
$ws.Range("J362").Value = 43
$ws.Range("K362").Value = 8000
$ws.Range("L362").Value = 9000
$ws.Range("M362").Value = 8512
$ws.Range("P362").Value = 2837

Write-Host "Done"
